$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = "sla_type"
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = "Minit"

$ws.Range("A44").Value = "sla_type"
$ws.Range("B44").Value = 2
$ws.Range("C44").Value = "Jam"

$ws.Range("A45").Value = "sla_type"
$ws.Range("B45").Value = 3
$ws.Range("C45").Value = "Hari"

$ws.Range("A46").Value = "severity"
$ws.Range("B46").Value = 1
$ws.Range("C46").Value = "Tidak Penting "

$ws.Range("A47").Value = "severity"
$ws.Range("B47").Value = 2
$ws.Range("C47").Value = "Kritikal"

$ws.Range("A48").Value = "severity"
$ws.Range("B48").Value = 3
$ws.Range("C48").Value = "Penting"

$ws.Range("A49").Value = "severity"
$ws.Range("B49").Value = 4
$ws.Range("C49").Value = "Sederhana"

$ws.Range("A50").Value = "severity"
$ws.Range("B50").Value = 5
$ws.Range("C50").Value = "Rendah"

$ws.PageSetup.Orientation = 1

$ws.Range("C50").Select()
